$d = $word.ActiveDocument

# --- Paragraph 1: split the Courier New run into 5 runs with proofErr markers ---
$p1 = $d.Paragraphs(1)
$s1 = $p1.Range.Start
$e1 = $p1.Range.End - 1
$r1 = $d.Range($s1, $e1)
$xmlP1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00313CEA"><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Please let us know if you have any discussion points you would like to highlight in your article and we may feature your article on our social media sites. Possible discussion points could highlight novel aspects of your paper that will be of interest to the engineering community, or perhaps raise questions regarding how other engineers view the matters discussed in your manuscript. Please note</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> discussion points should be no longer than a few sentences as they will be posted on the journal </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Facebook</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> page when your paper is published.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xmlP1)

# --- Paragraph 4 (the "AirStar" paragraph): split the run containing AirStar ---
$p4 = $d.Paragraphs(4)
$s4 = $p4.Range.Start
$e4 = $p4.Range.End - 1
$r4 = $d.Range($s4, $e4)
$xmlP4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">It is well know that tuning of multivariable adaptive control systems is a very time consuming task. It is often hard to decide whether the current set of </w:t></w:r><w:r w:rsidR="006D4E44"><w:t>control system parameters will</w:t></w:r><w:r><w:t xml:space="preserve"> guarantee </w:t></w:r><w:r w:rsidR="006D4E44"><w:t xml:space="preserve">the </w:t></w:r><w:r><w:t>best performance of a newly designed system i</w:t></w:r><w:r w:rsidR="006D4E44"><w:t>n</w:t></w:r><w:r><w:t xml:space="preserve"> real operational conditions. Thus, the paper presents an engineering approach to tuning of an L1 adaptive controller that is based on a highly uniform and very economical sampling of a multiple variable design space of </w:t></w:r><w:r w:rsidR="006D4E44"><w:t xml:space="preserve">the </w:t></w:r><w:r><w:t xml:space="preserve">desired control system parameters that enables construction of the Pareto front in the multidimensional control metrics space. The approach provided significant insight in the design and flight testing of an </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>AirStar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> jet airplane</w:t></w:r><w:r w:rsidR="006D4E44"><w:t xml:space="preserve"> implementing L1 adaptive autopilot</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r4.InsertXML($xmlP4)

# --- Append two new paragraphs after paragraph 4: one empty, one with new content ---
$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs(5)
$p5.Range.InsertParagraphAfter()

# --- Fill paragraph 6 with the new "Parameter Space Investigation" text ---
$p6 = $d.Paragraphs(6)
$s6 = $p6.Range.Start
$e6 = $p6.Range.End - 1
$r6 = $d.Range($s6, $e6)
$xmlP6 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>The use of the P</w:t></w:r><w:r><w:t xml:space="preserve">arameter  </w:t></w:r><w:r><w:t>S</w:t></w:r><w:r><w:t xml:space="preserve">pace </w:t></w:r><w:r><w:t>I</w:t></w:r><w:r><w:t>nvestigation</w:t></w:r><w:r><w:t xml:space="preserve"> method together with the MOVI software provided significant insights into the design and optimization of the L1 flight control system for the AirSTAR GTM aircraft, which was successfully flight tested by NASA.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r6.InsertXML($xmlP6)

Write-Output "done"
